# Add an "Execution Status" tracking block (columns H:K) to the WBS sheet
# and backfill automated status/date values for every WBS row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1, 8).Value  = "Execution Status"
$ws.Cells.Item(1, 9).Value  = "Status Updated On"
$ws.Cells.Item(1, 10).Value = "Started On"
$ws.Cells.Item(1, 11).Value = "Completed On"

# Column I holds a real date serial (2026-02-24) formatted as yyyy-mm-dd.
$statusUpdatedOn = 46077

# --- Row groups -------------------------------------------------------------
# Rows 2-15  : Done,        started/completed 2026-03-06
# Rows 16-26 : Done,        started/completed 2026-03-13
# Rows 27-49 : Partial,     started 2026-04-03 (not completed yet)
# Rows 50-137: Not Started, no start/completion date yet

for ($row = 2; $row -le 137; $row++) {

    if ($row -le 15) {
        $status = "Done"
        $startedOn = "2026-03-06"
        $completedOn = "2026-03-06"
    } elseif ($row -le 26) {
        $status = "Done"
        $startedOn = "2026-03-13"
        $completedOn = "2026-03-13"
    } elseif ($row -le 49) {
        $status = "Partial"
        $startedOn = "2026-04-03"
        $completedOn = $null
    } else {
        $status = "Not Started"
        $startedOn = $null
        $completedOn = $null
    }

    $ws.Cells.Item($row, 8).Value = $status

    $ws.Cells.Item($row, 9).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 9).Value = $statusUpdatedOn

    if ($startedOn) {
        $ws.Cells.Item($row, 10).NumberFormat = "yyyy-mm-dd"
        $ws.Cells.Item($row, 10).Value = "'" + $startedOn
    }

    if ($completedOn) {
        $ws.Cells.Item($row, 11).NumberFormat = "yyyy-mm-dd"
        $ws.Cells.Item($row, 11).Value = "'" + $completedOn
    }
}
